# Updated symbol list (cryptos.xlsx) - price/volume refresh + row reorder for rows 15-18
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '246.11'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '29.84'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '-0.96%'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.151'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '-0.65%'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.05753'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '0.17%'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '6.652'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '0.94%'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.239'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '6.54%'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.8488'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '-1.00%'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.8540'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '-2.23%'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1388'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '1.66%'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07086'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '0.30%'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.03256'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '11.45%'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.09374'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '-0.22%'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.001527'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '0.93%'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.005918'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '-2.68%'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.522'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '0.46%'
$ws.Range('B17').Value = 'BTSEToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.222'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '-2.06%'
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.01022'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '-0.28%'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.03360'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '2.64%'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.1315'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '0.60%'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.496'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '-2.74%'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.04131'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '-0.23%'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.1409'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.001228'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '1.09%'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.004144'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '-8.14%'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0001200'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '1.77%'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0001448'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.03748'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '-1.12%'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.1070'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '-0.17%'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.002299'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '4.61%'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.002949'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '-48.39%'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.009963'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '5.32%'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005524'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '8.41%'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '0.06%'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.07097'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.002467'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '-10.09%'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00002099'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '0.06%'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0001999'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '0.06%'
